# ArbeitRechercheProtokoll.xlsx — add two new work-log entries (rows 17 & 18)
# to the "Arbeitsprotokoll" sheet, describing the RGB-DEM / 3D-map scripting
# work, and update the sheet view to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsprotokoll")
$ws.Activate()

# --- new rows -----------------------------------------------------------
# Row 16 is the last existing entry; clone its formatting onto rows 17/18
# so the new cells pick up the same date / text styles used throughout the
# log, then overwrite the values.
$ws.Range("A16:B16").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A16:B16").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)   # xlPasteFormats

# Row 17's note is long / multi-paragraph like rows 6 & 10 -> reuse their
# word-wrap formatting for column B specifically.
$ws.Range("B10").Copy()
$ws.Range("B17").PasteSpecial(-4122)       # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Rows.Item(17).RowHeight = 64
$ws.Rows.Item(18).RowHeight = 64

$ws.Range("A17").Value = 45439
$ws.Range("B17").Value = "Bis hierhin auf versch. Wege mit Physiksim rumgepröbelt, nicht erfolgreich --> Rein analytisches Vorgehen`nHerunterladen von schweizer geodatene, zusammenhängen in eine grosse datei, reprojezieren von bessel auf wgs84 und RGB Format, darstellung in Visualisierung. Erstmals sind 3d- Landschaften in der Schweiz dargestellt"

$ws.Range("A18").Value = 45440
$ws.Range("B18").Value = "Herunterladen von weiteren Kantonen, konvertieren und erste versuche mit webhosting auf öffentlichem server"

# --- view state -----------------------------------------------------------
# Zoom out a bit and scroll down to the new entries, landing the selection
# just past the last row (matches where the author's cursor ended up).
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A19").Select() | Out-Null
